$d = $word.ActiveDocument

function Dump($label) {
    Write-Output "=== $label (count=$($d.Paragraphs.Count)) ==="
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        Write-Output "$i [$($p.Range.Start)-$($p.Range.End)] : [$($p.Range.Text)]"
    }
}

# =====================================================================
# Paragraph 1 (title): "Progress report February 9th 2018" (underlined)
#   -> "23" + "rd"(superscript) + " February 2018" (en-GB, no underline)
# Insert at the very start of the document (position 0): a lone <w:p>
# here becomes a genuine standalone new paragraph (position 0 is the
# only spot where this does not merge into neighboring content).
# =====================================================================
$xml1 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>23</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/><w:lang w:val="en-GB"/></w:rPr><w:t>rd</w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> February 2018</w:t></w:r></w:p>
'@
$d.Range(0, 0).InsertXML($xml1)

# Remove the old title paragraph (now paragraph 2), mark and all.
$oldTitle = $d.Paragraphs.Item(2)
$d.Range($oldTitle.Range.Start, $oldTitle.Range.End).Delete()

# =====================================================================
# Paragraph 2 (body): "This week we created the tables on JDBC, ..."
#   -> "We have finished the pojos and the logging scene of the GUI."
#      (with the _GoBack bookmark relocated inside it)
# Mid-document inserts need a real <w:p>...</w:p> followed by an empty
# dummy <w:p/> so the genuine paragraph break survives; the dummy's
# leftover empty paragraph is then deleted as its own discrete range
# (deleting it together with adjoining text in one range.Delete() call
# does not work reliably), then the old text paragraph is deleted.
# =====================================================================
$p2 = $d.Paragraphs.Item(2)
$xml2 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">We </w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">have finished the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>pojos</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> and the logging scene of </w:t></w:r><w:bookmarkStart w:id="100" w:name="_GoBack_NEW"/><w:bookmarkEnd w:id="100"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>the GUI.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
'@
$d.Range($p2.Range.Start, $p2.Range.Start).InsertXML($xml2)

# Delete the stray empty dummy paragraph (item 3) on its own.
$dummy2 = $d.Paragraphs.Item(3)
$d.Range($dummy2.Range.Start, $dummy2.Range.End).Delete()

# Delete the old body paragraph (now item 3) on its own.
$oldBody = $d.Paragraphs.Item(3)
$d.Range($oldBody.Range.Start, $oldBody.Range.End).Delete()

# =====================================================================
# Paragraph 3: "For next week, we will design the GUI ..."
#   -> "For next week we will do the controllers of all of them."
# =====================================================================
$p3 = $d.Paragraphs.Item(3)
$xml3 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>For next week we will do the controllers of all of them.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
'@
$d.Range($p3.Range.Start, $p3.Range.Start).InsertXML($xml3)

# Delete the stray empty dummy paragraph (item 4) on its own.
$dummy3 = $d.Paragraphs.Item(4)
$d.Range($dummy3.Range.Start, $dummy3.Range.End).Delete()

# Delete the old "For next week..." paragraph (now item 4). It still
# carries the original _GoBack bookmark, which disappears with it.
$oldP3 = $d.Paragraphs.Item(4)
$d.Range($oldP3.Range.Start, $oldP3.Range.End).Delete()

Dump "after all three paragraph rewrites"

# =====================================================================
# Drop the trailing empty paragraph (old paragraph 4) that followed
# "For next week..." in the source document. This paragraph sits at
# the very end of the document, so its range.Delete() must start one
# position earlier (inside the previous paragraph) to actually take
# effect.
# =====================================================================
$pCount = $d.Paragraphs.Count
$pLast = $d.Paragraphs.Item($pCount)
$pPrev = $d.Paragraphs.Item($pCount - 1)
$d.Range($pPrev.Range.End - 1, $pLast.Range.End).Delete()

Dump "after dropping trailing empty paragraph"

# =====================================================================
# Bookmark cleanup: rename our placeholder "_GoBack_NEW" bookmark back
# to "_GoBack" (the original "_GoBack" bookmark was already removed
# together with the old "For next week..." paragraph text above).
# (Bookmarks.Count / foreach-iteration are unreliable in this host, so
# look the bookmark up by name directly instead.)
# =====================================================================
$nb = $d.Bookmarks.Item("_GoBack_NEW")
$nbRange = $d.Range($nb.Start, $nb.End)
$nb.Delete()
$d.Bookmarks.Add("_GoBack", $nbRange)

Dump "final"
